$d = $word.ActiveDocument

# The document's single paragraph only carries the auto-generated
# "_GoBack" bookmark (Word drops this in every time you save after moving
# the cursor). Removing it collapses the paragraph to an empty <w:p/>,
# which is exactly what the target revision shows.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
